$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Change Area column (B) from "ITA17" to "SLO" for all data rows (2-97)
$ws.Range("B2:B97").Value = "SLO"

# Row 94 (before the delete) currently holds the CHLAGLA record that needs to be removed.
# Update rows 95-97 values "in place" is unnecessary because deleting row 94 shifts
# rows 95-97 up to 94-96 automatically, bringing the correct data with them.
# So we just need to delete the row containing SpecCode "CHLAGLA" (row 94).
$ws.Rows.Item(94).Delete()

